$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E columns to Text format so numeric-looking strings (e.g. "259.52", "98.826.48")
# are preserved as text instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "98.826.48"
$ws.Range("E2").Value = "  +1.90%  "
$ws.Range("D3").Value = "3.372.90"
$ws.Range("E3").Value = "  +8.36%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "259.52"
$ws.Range("E5").Value = "  +8.80%  "
$ws.Range("D6").Value = "630.17"
$ws.Range("E6").Value = "  +3.43%  "
$ws.Range("D7").Value = "1.42"
$ws.Range("E7").Value = "  +27.24%  "
$ws.Range("D8").Value = "0.394"
$ws.Range("E8").Value = "  +2.80%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").Value = "0.885"
$ws.Range("E10").Value = "  +11.97%  "
$ws.Range("D11").Value = "3.370.33"
$ws.Range("E11").Value = "  +8.50%  "
$ws.Range("D12").Value = "0.199"
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("D13").Value = "98.664.53"
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").Value = "36.22"
$ws.Range("E14").Value = "  +7.09%  "
$ws.Range("D15").Value = "0.0000249"
$ws.Range("E15").Value = "  +3.45%  "
$ws.Range("D16").Value = "3.984.60"
$ws.Range("E16").Value = "  +7.98%  "
$ws.Range("D17").Value = "5.54"
$ws.Range("E17").Value = "  +3.08%  "
$ws.Range("D18").Value = "3.368.68"
$ws.Range("E18").Value = "  +8.57%  "
$ws.Range("D19").Value = "3.59"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").Value = "15.30"
$ws.Range("E20").Value = "  +4.98%  "
$ws.Range("D21").Value = "495.34"
$ws.Range("E21").Value = "  -4.09%  "
$ws.Range("D22").Value = "6.21"
$ws.Range("E22").Value = "  +9.30%  "
$ws.Range("E23").Value = "  +9.24%  "
$ws.Range("D24").Value = "9.33"
$ws.Range("E24").Value = "  +5.68%  "
$ws.Range("D25").Value = "5.74"
$ws.Range("E25").Value = "  +4.16%  "
$ws.Range("D26").Value = "89.40"
$ws.Range("E26").Value = "  +2.94%  "
$ws.Range("D27").Value = "12.06"
$ws.Range("E27").Value = "  +3.75%  "
$ws.Range("D28").Value = "3.557.80"
$ws.Range("E28").Value = "  +8.50%  "
$ws.Range("D29").Value = "0.288"
$ws.Range("E29").Value = "  +20.56%  "
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").Value = "0.196"
$ws.Range("E31").Value = "  +12.35%  "
$ws.Range("D32").Value = "0.137"
$ws.Range("E32").Value = "  +8.57%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "9.65"
$ws.Range("E33").Value = "  +7.09%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").Value = "28.04"
$ws.Range("E35").Value = "  +5.37%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.152"
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "7.36"
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("D38").Value = "1.98"
$ws.Range("E38").Value = "  +5.63%  "
$ws.Range("D39").Value = "0.465"
$ws.Range("E39").Value = "  +6.53%  "
$ws.Range("D40").Value = "500.92"
$ws.Range("E40").Value = "  +3.31%  "
$ws.Range("D41").Value = "24.85"
$ws.Range("E41").Value = "  +2.69%  "
$ws.Range("B42").Value = "MantraDAO"
$ws.Range("C42").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D42").Value = "3.83"
$ws.Range("E42").Value = "  +4.76%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "1.28"
$ws.Range("E43").Value = "  +3.12%  "
$ws.Range("D44").Value = "3.36"
$ws.Range("E44").Value = "  +5.52%  "
$ws.Range("D45").Value = "0.787"
$ws.Range("E45").Value = "  +14.29%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "160.84"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").Value = "1.96"
$ws.Range("E48").Value = "  +3.44%  "
$ws.Range("D49").Value = "0.838"
$ws.Range("E49").Value = "  +15.09%  "
$ws.Range("D50").Value = "4.68"
$ws.Range("E50").Value = "  +6.87%  "
$ws.Range("D51").Value = "46.34"

# Restore default (unstyled) cell style on the Price/Volume columns so only
# the cell content changed and no stray number formatting remains applied.
$ws.Range("D2:E51").Style = "Normal"

